$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
  @{ r = 2; vals = @("ECs", "Hgf", "Cd44", "ECs", 2, 0.6666666666666666, 5.231719, 15.695157, 0.362499186434781, 0.362499186434781, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 1367.467774258957, 12307.20996833061, 0.2526139809281543, 0.2526139809281543) }
  @{ r = 3; vals = @("ECs", "Hgf", "Cd44", "FAPs", 2, 0.6666666666666666, 5.231719, 15.695157, 0.362499186434781, 0.362499186434781, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 167.412015456732, 1506.708139110588, 0.03092622471681166, 0.03092622471681166) }
  @{ r = 4; vals = @("ECs", "Hgf", "Cd44", "sCs", 2, 0.6666666666666666, 5.231719, 15.695157, 0.362499186434781, 0.362499186434781, 3, 1, 81.699019, 245.097057, 0.2178183668945166, 0.2178183668945167, 427.426309983661, 3846.836789852949, 0.0789589807898149, 0.07895898078981492) }
  @{ r = 5; vals = @("FAPs", "Hgf", "Cd44", "ECs", 3, 1, 8.316945333333335, 24.950836, 0.5762706133406404, 0.5762706133406403, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 2173.88485956657, 19564.96373609913, 0.4015843874289062, 0.4015843874289061) }
  @{ r = 6; vals = @("FAPs", "Hgf", "Cd44", "FAPs", 3, 1, 8.316945333333335, 24.950836, 0.5762706133406404, 0.5762706133406403, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 266.137493373936, 2395.237440365424, 0.04916390202457448, 0.04916390202457447) }
  @{ r = 7; vals = @("FAPs", "Hgf", "Cd44", "sCs", 3, 1, 8.316945333333335, 24.950836, 0.5762706133406404, 0.5762706133406403, 3, 1, 81.699019, 245.097057, 0.2178183668945166, 0.2178183668945167, 679.4862748099615, 6115.376473289652, 0.1255223238871598, 0.1255223238871598) }
  @{ r = 8; vals = @("sCs", "Hgf", "Cd44", "ECs", 3, 1, 0.8836963333333333, 2.651089, 0.06123020022457864, 0.06123020022457864, 3, 1, 261.380203, 784.1406089999999, 0.6968677182772199, 0.6968677182772199, 230.9807269970223, 2078.826542973201, 0.04266934992015944, 0.04266934992015944) }
  @{ r = 9; vals = @("sCs", "Hgf", "Cd44", "FAPs", 3, 1, 0.8836963333333333, 2.651089, 0.06123020022457864, 0.06123020022457864, 3, 1, 31.999428, 95.998284, 0.08531391482826334, 0.08531391482826335, 28.277777192364, 254.499994731276, 0.005223788086877213, 0.005223788086877214) }
  @{ r = 10; vals = @("sCs", "Hgf", "Cd44", "sCs", 3, 1, 0.8836963333333333, 2.651089, 0.06123020022457864, 0.06123020022457864, 3, 1, 81.699019, 245.097057, 0.2178183668945166, 0.2178183668945167, 72.19712352723033, 649.774111745073, 0.01333706221754198, 0.01333706221754199) }
)

foreach ($row in $rowData) {
  $rr = $row.r
  $vals = $row.vals
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item($rr, $col).Value = $vals[$i]
  }
}

Write-Output "Updated rows 2-10 with new Hgf-Cd44 sending/target cluster matrix"
